$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.252454666666667
$ws.Range("N2").Value = 24.757364
$ws.Range("O2").Value = 0.05349680956196952
$ws.Range("P2").Value = 0.05349680956196953
$ws.Range("Q2").Value = 10.86212840590667
$ws.Range("R2").Value = 97.75915565316001
$ws.Range("S2").Value = 0.05349680956196952
$ws.Range("T2").Value = 0.05349680956196953

# Row 3
$ws.Range("O3").Value = 0.5638948237978928
$ws.Range("P3").Value = 0.5638948237978929
$ws.Range("S3").Value = 0.5638948237978928
$ws.Range("T3").Value = 0.5638948237978929

# Row 4
$ws.Range("M4").Value = 57.81408433333333
$ws.Range("N4").Value = 173.442253
$ws.Range("O4").Value = 0.3747817085348802
$ws.Range("P4").Value = 0.3747817085348802
$ws.Range("Q4").Value = 76.09663222206333
$ws.Range("R4").Value = 684.8696899985699
$ws.Range("S4").Value = 0.3747817085348802
$ws.Range("T4").Value = 0.3747817085348802

# Row 5
$ws.Range("M5").Value = 1.207345666666667
$ws.Range("N5").Value = 3.622037
$ws.Range("O5").Value = 0.007826658105257385
$ws.Range("P5").Value = 0.007826658105257386
$ws.Range("Q5").Value = 1.589144586836667
$ws.Range("R5").Value = 14.30230128153
$ws.Range("S5").Value = 0.007826658105257385
$ws.Range("T5").Value = 0.007826658105257386
